# Refresh cryptos list figures (price + 1h volume change) to match the
# latest GitHub Actions data pull. A handful of Price cells look like
# plain numbers (e.g. "1.00", "0.380") but must stay literal text - the
# sheet stores every Price/Volume cell as text, so NumberFormat is forced
# to "@" before writing those specific cells to stop Excel from
# re-interpreting them as numbers and dropping the formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '89.844.37'
$ws.Range("E2").Value = '  +3.63%  '
$ws.Range("D3").Value = '3.207.73'
$ws.Range("E3").Value = '  +2.17%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '218.31'
$ws.Range("E5").Value = '  +6.90%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '630.21'
$ws.Range("E6").Value = '  +4.10%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.395'
$ws.Range("E7").Value = '  +8.20%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.698'
$ws.Range("E8").Value = '  +7.03%  '
$ws.Range("E9").Value = '  +0.11%  '
$ws.Range("D10").Value = '3.206.27'
$ws.Range("E10").Value = '  +2.31%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.582'
$ws.Range("E11").Value = '  +10.65%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.180'
$ws.Range("E12").Value = '  +2.68%  '
$ws.Range("E13").Value = '  +10.52%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '33.95'
$ws.Range("E14").Value = '  +7.04%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.43'
$ws.Range("E15").Value = '  +4.37%  '
$ws.Range("D16").Value = '3.809.89'
$ws.Range("E16").Value = '  +2.71%  '
$ws.Range("D17").Value = '89.772.64'
$ws.Range("E17").Value = '  +4.01%  '
$ws.Range("D18").Value = '3.230.63'
$ws.Range("E18").Value = '  +3.67%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0000234'
$ws.Range("E19").Value = '  +82.77%  '
$ws.Range("E20").Value = '  +17.59%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.65'
$ws.Range("E21").Value = '  +3.15%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '440.63'
$ws.Range("E22").Value = '  +7.79%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.72'
$ws.Range("E23").Value = '  +4.06%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.15'
$ws.Range("E24").Value = '  +2.44%  '
$ws.Range("B25").Value = 'NEARProtocol'
$ws.Range("C25").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.36'
$ws.Range("E25").Value = '  +5.24%  '
$ws.Range("B26").Value = 'Aptos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.18'
$ws.Range("E26").Value = '  +4.65%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '83.27'
$ws.Range("E27").Value = '  +14.75%  '
$ws.Range("D28").Value = '3.448.25'
$ws.Range("E28").Value = '  +4.54%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  +0.05%  '
$ws.Range("E30").Value = '  +0.87%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.00'
$ws.Range("E31").Value = '  +0.05%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.18'
$ws.Range("E32").Value = '  +41.73%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '8.57'
$ws.Range("E33").Value = '  +4.68%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '551.47'
$ws.Range("E34").Value = '  +3.77%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '7.15'
$ws.Range("E35").Value = '  +10.35%  '
$ws.Range("B36").Value = 'PancakeSwap'
$ws.Range("C36").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.93'
$ws.Range("E36").Value = '  +5.05%  '
$ws.Range("B37").Value = 'Fetch.AI'
$ws.Range("C37").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.34'
$ws.Range("E37").Value = '  +5.34%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '22.57'
$ws.Range("E38").Value = '  +5.12%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '22.40'
$ws.Range("E39").Value = '  +3.00%  '
$ws.Range("E40").Value = '  -1.29%  '
$ws.Range("E41").Value = '  +0.03%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.96'
$ws.Range("E42").Value = '  +4.33%  '
$ws.Range("B43").Value = 'USDe'
$ws.Range("C43").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.00'
$ws.Range("E43").Value = '  -0.04%  '
$ws.Range("B44").Value = 'PolygonEcosystemToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.380'
$ws.Range("E44").Value = '  +3.87%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '147.21'
$ws.Range("E45").Value = '  -1.05%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '174.99'
$ws.Range("E46").Value = '  +2.88%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '43.89'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.777'
$ws.Range("E48").Value = '  +14.22%  '
$ws.Range("E49").Value = '  +0.12%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.26'
$ws.Range("E50").Value = '  +2.14%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.628'
$ws.Range("E51").Value = '  +8.43%  '
